$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the (hidden) "_GoBack" bookmark from the "select * from customer c"
#    paragraph. Word hides bookmarks whose name starts with "_" from the
#    normal Bookmarks collection enumeration, but they can still be reached
#    (and deleted) by name.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
}

# ---------------------------------------------------------------------------
# 2) Append new paragraphs at the end of the document (after the last
#    paragraph, which holds the second screenshot), right before the
#    section break:
#       <empty>
#       <empty>
#       <empty>
#       <empty>
#       https://github.com/QuangDuc1512/KTPMUD_DV19_CK.git   (+ "_GoBack")
#       <empty>
# ---------------------------------------------------------------------------

# Create one fresh paragraph after the current last paragraph so that the
# insertion point below is no longer sitting exactly at the end of the
# story (inserting XML right at end-of-story would clobber the previous
# paragraph's content instead of appending after it).
$tailRange = $d.Paragraphs.Last.Range
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()

# Insertion point: start of that brand-new trailing paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $lastPara.Range.Duplicate
$insertPoint.Collapse(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$url = "https://github.com/QuangDuc1512/KTPMUD_DV19_CK.git"

$xmlFrag = "<w:p $wNs/>" +
           "<w:p $wNs/>" +
           "<w:p $wNs/>" +
           "<w:p $wNs/>" +
           "<w:p $wNs><w:r><w:t>$url</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>" +
           "<w:p $wNs/>"

$insertPoint.InsertXML($xmlFrag)
